# Add a new wishlist item ("Papusa") as row 26 of the Craciun sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Write the new row's cells. The order below controls the order in which
# new entries are appended to the shared-strings table, so it is chosen to
# reproduce the expected sharedStrings.xml ordering:
#   95 = Name ("Papusa "), 96 = Link (Amazon URL), 97 = Price ("23 EUR"),
#   98 = Image (amazon media URL)
$ws.Range("A26").Value = "Papusa "
$ws.Range("C26").Value = "https://www.amazon.de/-/en/105733519-Princess-Movable-Suitable-Children/dp/B08QYTD24X/ref=sr_1_2?crid=K5EE4ZX6PC2F&dib=eyJ2IjoiMSJ9.6WxTQND8ZRjVlveNJY2qRw0Mm4LH90lyvKyzz7nDafVgwRCnIYzN0bWBiMWZiMP8P163-QVHZGnO9vDq3RJAZwEAzmDQXfBGsgqWPa4-0g7KruRsNtj8nfuKPUBWSzuOZ-3hLnN7hABHEDEYMhLbj_OGraTMT0BRQWcsBS1khICmFnPU3eHBo41bmEBw0YXRe8BbNMYczPvSacihlWG3O1IyUVgp8i5f3cKO4s39KRVhD4gPPFJEUYL93RC4pV6NOYax081G7OZyfV0wcViHLvD-wO1o6ghtjUnqc7IxnQk.KN4nYMyUPUc1xxm5YUWrMV7Lym367vaok8RY7lT9OrI&dib_tag=se&keywords=steffi+love+pferd&qid=1729591110&sprefix=stefi+lo%2Caps%2C116&sr=8-2#customerReviews"
$ws.Range("D26").Value = "23 EUR"
$ws.Range("B26").Value = "https://m.media-amazon.com/images/I/51WGSlK3v6L._AC_SX679_.jpg"

# Match the selection shown in the author's saved view.
$ws.Range("B26").Select()
